$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price/Volume columns so numeric-looking strings
# (e.g. "64.273.96", "0.999") are preserved exactly as text, matching
# the original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.273.96"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.490.29"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D5").Value = "587.51"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "134.29"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.487"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "7.27"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "4.085.73"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.120"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "3.490.35"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "25.76"
$ws.Range("E16").Value = "  -7.19%  "
$ws.Range("D17").Value = "64.339.64"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "9.87"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "13.61"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").Value = "394.37"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "3.629.17"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").Value = "74.75"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "5.73"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  -5.82%  "
$ws.Range("D32").Value = "8.24"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "3.510.63"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "23.41"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "5.14"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").Value = "166.10"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").Value = "0.0781"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").Value = "0.805"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "25.15"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "1.16"
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("D48").Value = "2.455.45"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "0.893"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "0.0260"
$ws.Range("E51").Value = "  -1.17%  "
